$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (strike) values per row, replacing the previous Strike# values
$kValues = @{
    2  = 5
    3  = 6
    4  = 2
    5  = 4
    6  = 3
    7  = 7
    8  = 3
    9  = 3
    10 = 3
    11 = 3
    12 = 3
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
